# SFC Country Group workbook update
# Add a new "Yemen" entry under the "Middle East - others (2)" SFC Country Group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new row. Column order (C, then A, then B) is chosen so that the
# shared-string table gets the three new strings appended in the same order
# they appear in the target workbook: "Middle East - others (2)", "YE", "Yemen".
$ws.Range("C37").Value = "Middle East - others (2)"
$ws.Range("A37").Value = "YE"
$ws.Range("B37").Value = "Yemen"

# Move the selection past the newly entered data, mirroring where the author
# left the cursor after typing the new row.
[void]$ws.Range("A38").Select()
